$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $ws.Name
Write-Host $wb.Worksheets.Count
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
